$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-4.36%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-7.14%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.107"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.05%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07752"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-5.37%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.255"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.00%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.626"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-10.63%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8810"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.58%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1025"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.07%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1745"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.88%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08956"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.67%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04427"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.60%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.23%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001264"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.44%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005803"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.33%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.31%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.436"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.18%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3281"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.84%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.008"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-5.38%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1340"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.40%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2787"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "11.79%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04180"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.38%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.58%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004085"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-6.11%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.33%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.05%"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-14.72%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05211"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.93%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.05%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1325"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-5.19%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006342"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.15%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001964"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.13%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008764"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "15.96%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3345"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.11%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006537"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-6.31%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "98.41%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002730"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-21.81%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
